$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: update title and link
$ws.Range("D9").Value = "10년차 고등학교 수학 강사의 Data Science 도전?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/math-teacher-for-data-science/#utm_source=rss&utm_medium=rss&utm_campaign=math-teacher-for-data-science"

# Row 26: update title
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 52: update title
$ws.Range("D52").Value = "Relative Risk Regression"
